$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared string label (country/variable name) for row 2
$ws.Range("A2").Value = "Sprudent_npv_over_gdp_gcs_adj"

# GCP max_gain, prudent map update -- refreshed numeric values for row 2
$ws.Range("B2").Value = 0.0151043838170846
$ws.Range("C2").Value = 0.00359379866853409
$ws.Range("F2").Value = 0.000410639911282762
$ws.Range("I2").Value = -0.00906753310710586
$ws.Range("K2").Value = 0.0197357646957889
$ws.Range("L2").Value = -0.0110574671035766
$ws.Range("M2").Value = 0.00281927431191736
$ws.Range("N2").Value = 0.0123989249979428
$ws.Range("O2").Value = 0.0070820295521727
$ws.Range("P2").Value = -0.00465620993605702
$ws.Range("R2").Value = 0.000289105637803429
$ws.Range("U2").Value = 0.000876068697373236
$ws.Range("V2").Value = 0.000623788839979223
$ws.Range("W2").Value = 0.000514636877426194
$ws.Range("X2").Value = 0.000805313115135587
$ws.Range("AB2").Value = 0.0273964563325208
$ws.Range("AD2").ClearContents()
$ws.Range("AG2").Value = 0.000156350839541344
$ws.Range("AH2").Value = 0.00181616904949319
$ws.Range("AI2").Value = 0.0102379030134157
$ws.Range("AJ2").Value = 0.00170710968242554
$ws.Range("AK2").Value = 0.000512474295276406
$ws.Range("AL2").Value = 0.00496401053353295
$ws.Range("AM2").Value = 0.00004134735942103
$ws.Range("AN2").Value = 0.00046950652686288
$ws.Range("AO2").Value = 0.000911133887080572
$ws.Range("AQ2").Value = -0.017371844776287
$ws.Range("AR2").Value = -0.0123825197143604
$ws.Range("AS2").Value = 0.000574858850398968
$ws.Range("AT2").Value = -0.00671104301656217
$ws.Range("AU2").Value = 0.000511633593603243
$ws.Range("AW2").Value = 0.000580556903979908
$ws.Range("AY2").Value = 0.014770077126678
$ws.Range("AZ2").Value = -0.000880917648465152
$ws.Range("BA2").Value = -0.0311435299324767
$ws.Range("BB2").Value = 0.0129126162029414
$ws.Range("BC2").Value = -0.00970686304218662
$ws.Range("BD2").Value = 0.00151572673886781
$ws.Range("BE2").Value = -0.00515433738859531
$ws.Range("BG2").ClearContents()
$ws.Range("BI2").Value = 0.00191845711714543
$ws.Range("BJ2").Value = 0.00101119869884897
$ws.Range("BK2").Value = 0.00867115959277107
$ws.Range("BL2").Value = 0.00768467925022399
$ws.Range("BM2").Value = 0.000435977988835038
$ws.Range("BN2").Value = -0.000062918905487202
$ws.Range("BO2").Value = 0.000997688315160532
$ws.Range("BP2").Value = 0.000627857482623868
$ws.Range("BR2").Value = 0.00193630465305304
$ws.Range("BT2").Value = 0.00911271981117367
$ws.Range("BU2").Value = -0.0045623288140523
$ws.Range("BW2").Value = 0.0000481930823676458
$ws.Range("BX2").Value = -0.00713749142471962
$ws.Range("CC2").Value = -0.00407972585902377
$ws.Range("CD2").Value = 0.00145364858923396
$ws.Range("CF2").ClearContents()
$ws.Range("CH2").Value = 0.00547561253493542
$ws.Range("CJ2").Value = 0.00810508653352903
$ws.Range("CM2").Value = 0.00563660850327913
$ws.Range("CO2").Value = 0.00688455444802085
$ws.Range("CQ2").Value = 0.00031973800393597
$ws.Range("CR2").Value = 0
$ws.Range("CS2").Value = -0.00594978512968433
$ws.Range("CT2").Value = -0.0113123376260774
$ws.Range("CU2").Value = -0.00285742455829926
$ws.Range("CX2").Value = 0.0115300983234055
$ws.Range("CY2").Value = 0
$ws.Range("CZ2").Value = 0.000358727919134943
$ws.Range("DB2").Value = 0.0056792195898242
$ws.Range("DC2").Value = -0.00521040600490127
$ws.Range("DD2").Value = 0.00804082587520934
$ws.Range("DG2").Value = 0.00531935055390146
$ws.Range("DH2").Value = 0.000722589352761333
$ws.Range("DJ2").Value = 0.0206837377157536
$ws.Range("DM2").Value = 0.0163720718627636
$ws.Range("DN2").Value = 0.000995848385274252
$ws.Range("DO2").Value = 0.000931642352597824
$ws.Range("DP2").Value = -0.0122990023591668
$ws.Range("DQ2").ClearContents()
$ws.Range("DR2").Value = 0.0118578424122401
$ws.Range("DS2").ClearContents()
$ws.Range("DU2").Value = 0.00423373735789115
$ws.Range("DV2").Value = -0.00205414744280956
$ws.Range("DW2").Value = 0.000408714782782181
$ws.Range("DX2").Value = 0.00321012755637867
$ws.Range("DY2").Value = 0.00485096659966935
$ws.Range("DZ2").Value = -0.00997287707281265
$ws.Range("EA2").Value = -0.000218165248823758
$ws.Range("EB2").Value = 0.000551340277322914
$ws.Range("ED2").Value = -0.00347339128967313
$ws.Range("EF2").Value = 0.0136249242800312
$ws.Range("EH2").Value = 0.0022793872973777
$ws.Range("EI2").Value = 0.00158428263815852
$ws.Range("EK2").Value = 0.0148941046820604
$ws.Range("EL2").Value = 0.0130299989590248
$ws.Range("EM2").Value = 0.001714630929229
$ws.Range("EN2").Value = 0.00631057285234826
$ws.Range("EP2").Value = 0.0112928448214244
$ws.Range("EQ2").Value = 0.000306592455004204
$ws.Range("ER2").Value = -0.0103121841994999
$ws.Range("ES2").Value = -0.00523376878821256
$ws.Range("ET2").Value = -0.00461662233276537
$ws.Range("EU2").Value = 0.000289396351640338
$ws.Range("EW2").Value = 0.0156906420238651
$ws.Range("EX2").Value = 0.00427737927459854
$ws.Range("EZ2").Value = 0.00543701416122187
$ws.Range("FB2").Value = 0.00244177260296653
$ws.Range("FG2").Value = 0.0064762855540297
$ws.Range("FH2").Value = 0.0104949467258887
$ws.Range("FJ2").Value = -0.000536334110083664
$ws.Range("FM2").Value = 0.000831871945226889
$ws.Range("FO2").Value = 0.0100872572897246
$ws.Range("FP2").Value = 0.00360472029934705
$ws.Range("FQ2").Value = 0.000112782239843026
$ws.Range("FS2").Value = 0.00318834924369617
$ws.Range("FT2").Value = 0.0000349231624274589
